# Update "想去人数" (want-to-go count) figures in the F column across sheets.
# Sheet "展览" (Exhibition)      -> xl/worksheets/sheet1.xml
# Sheet "演出" (Performance)     -> xl/worksheets/sheet2.xml
# Sheet "全部类型" (All types)    -> xl/worksheets/sheet4.xml

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsPerformance = $wb.Worksheets.Item("演出")
$wsAllTypes = $wb.Worksheets.Item("全部类型")

# --- 展览 (Exhibition) sheet ---
$wsExhibition.Range("F2").Value = 212
$wsExhibition.Range("F3").Value = 120
$wsExhibition.Range("F5").Value = 999
$wsExhibition.Range("F6").Value = 5571
$wsExhibition.Range("F7").Value = 502
$wsExhibition.Range("F8").Value = 704
$wsExhibition.Range("F9").Value = 962
$wsExhibition.Range("F10").Value = 827
$wsExhibition.Range("F13").Value = 590
$wsExhibition.Range("F14").Value = 31
$wsExhibition.Range("F17").Value = 1866
$wsExhibition.Range("F18").Value = 1477
$wsExhibition.Range("F19").Value = 935
$wsExhibition.Range("F21").Value = 197
$wsExhibition.Range("F22").Value = 341
$wsExhibition.Range("F23").Value = 559
$wsExhibition.Range("F25").Value = 1056
$wsExhibition.Range("F28").Value = 2985
$wsExhibition.Range("F30").Value = 104
$wsExhibition.Range("F31").Value = 65
$wsExhibition.Range("F32").Value = 125
$wsExhibition.Range("F33").Value = 38
$wsExhibition.Range("F34").Value = 399
$wsExhibition.Range("F35").Value = 14
$wsExhibition.Range("F39").Value = 296
$wsExhibition.Range("F40").Value = 734
$wsExhibition.Range("F41").Value = 92
$wsExhibition.Range("F42").Value = 53
$wsExhibition.Range("F43").Value = 59

# --- 演出 (Performance) sheet ---
$wsPerformance.Range("F4").Value = 200
$wsPerformance.Range("F6").Value = 138

# --- 全部类型 (All types) sheet ---
$wsAllTypes.Range("F3").Value = 212
$wsAllTypes.Range("F4").Value = 120
$wsAllTypes.Range("F5").Value = 999
$wsAllTypes.Range("F7").Value = 5571
$wsAllTypes.Range("F8").Value = 502
$wsAllTypes.Range("F9").Value = 704
$wsAllTypes.Range("F11").Value = 200
$wsAllTypes.Range("F12").Value = 962
$wsAllTypes.Range("F13").Value = 827
$wsAllTypes.Range("F15").Value = 138
$wsAllTypes.Range("F18").Value = 590
$wsAllTypes.Range("F19").Value = 31
$wsAllTypes.Range("F23").Value = 1866
$wsAllTypes.Range("F24").Value = 1477
$wsAllTypes.Range("F25").Value = 935
$wsAllTypes.Range("F26").Value = 197
$wsAllTypes.Range("F27").Value = 341
$wsAllTypes.Range("F29").Value = 559
$wsAllTypes.Range("F31").Value = 1056
$wsAllTypes.Range("F32").Value = 2985
$wsAllTypes.Range("F34").Value = 104
$wsAllTypes.Range("F35").Value = 65
$wsAllTypes.Range("F36").Value = 125
$wsAllTypes.Range("F37").Value = 38
$wsAllTypes.Range("F38").Value = 399
$wsAllTypes.Range("F39").Value = 14
$wsAllTypes.Range("F42").Value = 296
$wsAllTypes.Range("F43").Value = 734
$wsAllTypes.Range("F44").Value = 92
$wsAllTypes.Range("F45").Value = 59
